$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.235000000000001
$ws.Range("C3").Value = -11.817
$ws.Range("A4").Value = -21.562
$ws.Range("B4").Value = 6.563
$ws.Range("C4").Value = -12.325
$ws.Range("E4").Value = 12.883
$ws.Range("B5").Value = 6.254
$ws.Range("A6").Value = -21.109
$ws.Range("A7").Value = -21.303
$ws.Range("B8").Value = 6.361
$ws.Range("C9").Value = -11.775
$ws.Range("C11").Value = -12.642
$ws.Range("E12").Value = 12.959
$ws.Range("C14").Value = -11.784
$ws.Range("A16").Value = -21.238
$ws.Range("B16").Value = 6.436
$ws.Range("E17").Value = 13.035
$ws.Range("C18").Value = -12.412
$ws.Range("A20").Value = -22.36
$ws.Range("D20").Value = -8.263999999999999
$ws.Range("E20").Value = 13.154
$ws.Range("B22").Value = 6.693000000000001
$ws.Range("C25").Value = -12.347
$ws.Range("E25").Value = 13.09
